# chore: Simplifies instructor information.
#
# All seeded instructors previously had distinct per-user passwords
# (Tutor!1991 / Tutor!1999 / Tutor!1998 / Tutor!1995). This replaces every
# instructor's password with a single shared generic value, "tutor".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Password column (D) for every instructor data row (2-5).
$ws.Range("D2:D5").Value = "tutor"

# Matches the author's last selection in the saved workbook.
$null = $ws.Range("D6").Select()
